$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Cap nhat phan trang" -- reset the paging columns (I and J) to 0 for
# every data row (rows 2 through 89).
for ($r = 2; $r -le 89; $r++) {
    $ws.Cells.Item($r, 9).Value = 0   # column I
    $ws.Cells.Item($r, 10).Value = 0  # column J
}
